$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Julio de 2020 a las 19:40"

# --- Swap country labels caused by re-sorting (values tied/crossed in rank) ---
# Ghana overtook Irlanda: row 57 becomes Ghana (with refreshed numbers),
# row 58 becomes Irlanda (keeping Irlanda's previous numbers).
$ws.Range("A57").Value = "Ghana"
$ws.Range("A58").Value = "Irlanda"

# Islas Malvinas / Groenlandia are tied, just swap label order.
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

# --- Numeric updates per country row (B=Casos totales, C=Nuevos casos, ---
# --- D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 3649049
$ws.Range("C4").Value = 32222
$ws.Range("D4").Value = 1647761
$ws.Range("E4").Value = 1860762
$ws.Range("G4").Value = 382
$ws.Range("H4").Value = 140526

# Row 6: India
$ws.Range("B6").Value = 1004383
$ws.Range("C6").Value = 34214
$ws.Range("D6").Value = 636541
$ws.Range("E6").Value = 342237
$ws.Range("G6").Value = 676
$ws.Range("H6").Value = 25605

# Row 9: Chile
$ws.Range("B9").Value = 323698
$ws.Range("C9").Value = 2493
$ws.Range("D9").Value = 295301
$ws.Range("E9").Value = 21107
$ws.Range("G9").Value = 104
$ws.Range("H9").Value = 7290

# Row 18: Turquia
$ws.Range("B18").Value = 216873
$ws.Range("C18").Value = 933
$ws.Range("D18").Value = 198820
$ws.Range("E18").Value = 12613
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = 5440

# Row 24: Canada
$ws.Range("B24").Value = 109080
$ws.Range("C24").Value = 251
$ws.Range("D24").Value = 72782
$ws.Range("E24").Value = 27473
$ws.Range("G24").Value = 15
$ws.Range("H24").Value = 8825

# Row 46: Israel
$ws.Range("B46").Value = 45607
$ws.Range("C46").Value = 1419
$ws.Range("D46").Value = 20268
$ws.Range("E46").Value = 24956
$ws.Range("G46").Value = 7
$ws.Range("H46").Value = 383

# Row 50: Barein
$ws.Range("E50").Value = 4119
$ws.Range("G50").Value = 4
$ws.Range("H50").Value = 121

# Row 57: now Ghana (refreshed numbers)
$ws.Range("B57").Value = 26125
$ws.Range("C57").Value = 695
$ws.Range("D57").Value = 22270
$ws.Range("E57").Value = 3716
$ws.Range("H57").Value = 139

# Row 58: now Irlanda (carries old Irlanda numbers down)
$ws.Range("B58").Value = 25683
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 23364
$ws.Range("E58").Value = 571
$ws.Range("H58").Value = 1748

# Row 60: Argelia
$ws.Range("B60").Value = 21355
$ws.Range("C60").Value = 585
$ws.Range("D60").Value = 15107
$ws.Range("E60").Value = 5196
$ws.Range("G60").Value = 12
$ws.Range("H60").Value = 1052

# Row 87: Estado de Palestina
$ws.Range("D87").Value = 1313
$ws.Range("E87").Value = 6050
$ws.Range("G87").Value = 5
$ws.Range("H87").Value = 49

# Row 112: Libano
$ws.Range("B112").Value = 2599
$ws.Range("C112").Value = 57
$ws.Range("D112").Value = 1485
$ws.Range("E112").Value = 1074
$ws.Range("G112").Value = 2
$ws.Range("H112").Value = 40

# Row 127: Libia
$ws.Range("B127").Value = 1652
$ws.Range("C127").Value = 63
$ws.Range("D127").Value = 379
$ws.Range("E127").Value = 1227
$ws.Range("G127").Value = 3
$ws.Range("H127").Value = 46
